$wb = $excel.ActiveWorkbook

# --- Update the conversion note text on sheet "Hoja1" (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.86 = 19271.71 pesos`n✅ 19271.71 pesos = 4.85 = 953.89 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 205.55
$wsTasas.Range("O10").Value = 3961.3
$wsTasas.Range("N12").Value = 3975
$wsTasas.Range("O12").Value = 196.75
